# Apply "finished integration frequency but bugs remain with some segs having freq na"
# edit to segments.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the worksheet (this also updates the _xlnm._FilterDatabase defined
#    name, which references the sheet by name).
# ---------------------------------------------------------------------------
$ws.Name = "segments.txt"

# ---------------------------------------------------------------------------
# 2) New row 62: segment "w", a consonant with freq=2483, pc=82.
#    A62 is set first so the new shared string "w" is registered before the
#    "ə:" string introduced below (matches original authoring order).
# ---------------------------------------------------------------------------
$ws.Range("A62").Value = "w"
$ws.Range("B62").Value = "consonant"
$ws.Range("B62").Style = "Normal"   # column B carries a column style; reset so B62 has no explicit style, like the rest of the new row
$ws.Range("C62").Value = "consonant"
$ws.Range("D62").Value = "consonant"
$ws.Range("E62").Value = "consonant"
$ws.Range("F62").Value = 2483
$ws.Range("G62").Value = 82

# ---------------------------------------------------------------------------
# 3) Row 15 (segment "êê"/ə): correct the phono value to the long vowel "ə:"
# ---------------------------------------------------------------------------
$ws.Range("B15").Value = "ə:"

# ---------------------------------------------------------------------------
# 4) Row 16 (segment "ââ"): this segment now has real frequency data filled
#    in, replacing the previous "NA" placeholders.
# ---------------------------------------------------------------------------
$ws.Range("E16").Value = "yes"

$ws.Range("F16").Value = 104
$ws.Range("F16").Font.Name = "Helvetica Neue"
$ws.Range("F16").Font.Size = 14
$ws.Range("F16").Font.Color = 3355443

$ws.Range("G16").Value = 3
$ws.Range("H16").Value = $null

$ws.Rows.Item(16).RowHeight = 18

# ---------------------------------------------------------------------------
# 5) Row 61 (segment "uu"/u:): in_stim corrected from "yes" to "consonant"
#    (data entry fix noted in the commit message).
# ---------------------------------------------------------------------------
$ws.Range("E61").Value = "consonant"

# ---------------------------------------------------------------------------
# 6) Update the view state: active selection moved to F17, scrolled so row 8
#    is at the top.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("F17").Select() | Out-Null

# Window geometry, as recorded in the commit.
$win.Left = 2540
$win.Top = 460
$win.Width = 20640
$win.Height = 14180
